$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A4").NumberFormat = "@"

$ws.Range("A2").Value = "337829999"
$ws.Range("B2").Value = "Tran Van Thanh"
$ws.Range("C2").Value = "nbk-vl"

$ws.Range("A3").Value = "285656899"
$ws.Range("B3").Value = "Le Cong Huy"
$ws.Range("C3").Value = "nbk-vl"

$ws.Range("A4").Value = "334442222"
$ws.Range("B4").Value = "Nguyen Thi Tho"
$ws.Range("C4").Value = "nbk-qn"

$ws.Range("A2:C4").ClearFormats()

$ws.Columns.Item(1).ColumnWidth = 10.2
$ws.Columns.Item(2).ColumnWidth = 14.3
